$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text interpretation (D/E columns hold text like "37.224.76" or "  +0.33%  ")
# while value-setting, then clear the temporary number-format so no residual
# per-cell style attribute is left behind (matches original cells, which carry no "s").
$fmtRange = $ws.Range("D2:E51")
$fmtRange.NumberFormat = "@"

$ws.Range("D2").Value = "37.243.79"
$ws.Range("E2").Value = "  +0.24%  "
$ws.Range("D3").Value = "2.077.03"
$ws.Range("E3").Value = "  -0.61%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "251.44"
$ws.Range("D6").Value = "0.675"
$ws.Range("E6").Value = "  +3.76%  "
$ws.Range("D7").Value = "62.99"
$ws.Range("E7").Value = "  +26.08%  "
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  -0.19%  "
$ws.Range("D9").Value = "61.88"
$ws.Range("E9").Value = "  +2.24%  "
$ws.Range("D10").Value = "0.387"
$ws.Range("E10").Value = "  +5.36%  "
$ws.Range("D11").Value = "0.0808"
$ws.Range("E11").Value = "  +9.50%  "
$ws.Range("E12").Value = "  +2.67%  "
$ws.Range("D13").Value = "15.69"
$ws.Range("E13").Value = "  +3.54%  "
$ws.Range("D14").Value = "2.377.20"
$ws.Range("E14").Value = "  -0.57%  "
$ws.Range("D15").Value = "0.828"
$ws.Range("E15").Value = "  -0.02%  "
$ws.Range("D16").Value = "5.43"
$ws.Range("E16").Value = "  +7.53%  "
$ws.Range("D17").Value = "2.081.63"
$ws.Range("E17").Value = "  +0.05%  "
$ws.Range("D18").Value = "37.185.79"
$ws.Range("E18").Value = "  +0.39%  "
$ws.Range("D19").Value = "74.94"
$ws.Range("E19").Value = "  +4.18%  "
$ws.Range("D20").Value = "0.0₃0928"
$ws.Range("E20").Value = "  +13.27%  "
$ws.Range("D21").Value = "15.15"
$ws.Range("E21").Value = "  +14.98%  "
$ws.Range("D22").Value = "5.48"
$ws.Range("E22").Value = "  +5.91%  "
$ws.Range("D23").Value = "240.30"
$ws.Range("E23").Value = "  +0.86%  "
$ws.Range("D24").Value = "1.00"
$ws.Range("E24").Value = "  -0.13%  "
$ws.Range("E25").Value = "  -0.75%  "
$ws.Range("D26").Value = "171.68"
$ws.Range("E26").Value = "  +1.75%  "
$ws.Range("D27").Value = "9.29"
$ws.Range("E27").Value = "  +0.10%  "
$ws.Range("D28").Value = "2.07"
$ws.Range("E28").Value = "  +3.76%  "
$ws.Range("D29").Value = "20.48"
$ws.Range("E29").Value = "  -1.22%  "
$ws.Range("D30").Value = "0.126"
$ws.Range("E30").Value = "  +3.48%  "
$ws.Range("D31").Value = "4.78"
$ws.Range("E31").Value = "  +7.10%  "
$ws.Range("D32").Value = "1.10"
$ws.Range("E32").Value = "  +3.07%  "
$ws.Range("D33").Value = "0.0639"
$ws.Range("E33").Value = "  +5.92%  "
$ws.Range("D34").Value = "4.45"
$ws.Range("E34").Value = "  +10.18%  "
$ws.Range("E35").Value = "  -1.01%  "
$ws.Range("D36").Value = "0.999"
$ws.Range("E36").Value = "  -0.20%  "
$ws.Range("D37").Value = "2.32"
$ws.Range("E37").Value = "  +2.99%  "
$ws.Range("E38").Value = "  -3.54%  "
$ws.Range("E39").Value = "  +23.75%  "
$ws.Range("E40").Value = "  +3.06%  "
$ws.Range("D41").Value = "18.87"
$ws.Range("E41").Value = "  +8.08%  "
$ws.Range("E42").Value = "  +2.55%  "
$ws.Range("E43").Value = "  +2.12%  "
$ws.Range("D44").Value = "99.06"
$ws.Range("E44").Value = "  +1.54%  "
$ws.Range("E45").Value = "  +24.24%  "
$ws.Range("E46").Value = "  +1.44%  "
$ws.Range("D47").Value = "2.58"
$ws.Range("E47").Value = "  +14.82%  "
$ws.Range("D48").Value = "4.56"
$ws.Range("E48").Value = "  +15.41%  "
$ws.Range("D49").Value = "1.311.51"
$ws.Range("E49").Value = "  +0.67%  "
$ws.Range("E50").Value = "  -0.73%  "
$ws.Range("D51").Value = "6.93"
$ws.Range("E51").Value = "  +1.12%  "

$fmtRange.ClearFormats()
